$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 963.2308
$ws.Cells.Item(38, 10).Value = 5399.5
$ws.Cells.Item(38, 12).Value = 16198.5
$ws.Cells.Item(38, 14).Value = -16942.5
$ws.Cells.Item(43, 8).Value = 2900
$ws.Cells.Item(43, 9).Value = 3325
$ws.Cells.Item(43, 11).Value = 3325
$ws.Cells.Item(43, 13).Value = -3256
$ws.Cells.Item(70, 8).Value = 101220950
$ws.Cells.Item(70, 9).Value = 337399200
$ws.Cells.Item(70, 10).Value = 1700
$ws.Cells.Item(70, 11).Value = 1012197600
$ws.Cells.Item(70, 12).Value = 5100
$ws.Cells.Item(70, 13).Value = -1012197330
$ws.Cells.Item(70, 14).Value = -5640
$ws.Cells.Item(73, 8).Value = 101220950
$ws.Cells.Item(73, 9).Value = 337399200
$ws.Cells.Item(73, 10).Value = 1700
$ws.Cells.Item(73, 11).Value = 1012197600
$ws.Cells.Item(73, 12).Value = 5100
$ws.Cells.Item(73, 13).Value = -1012196664
$ws.Cells.Item(73, 14).Value = -6972
$ws.Cells.Item(86, 8).Value = 8919.583000000001
$ws.Cells.Item(86, 9).Value = 3708
$ws.Cells.Item(86, 10).Value = 12642.143
$ws.Cells.Item(86, 11).Value = 3708
$ws.Cells.Item(86, 12).Value = 12642.143
$ws.Cells.Item(86, 13).Value = -2585
$ws.Cells.Item(86, 14).Value = -14888.143
$ws.Cells.Item(88, 8).Value = 2438.3125
$ws.Cells.Item(88, 9).Value = 2667.5
$ws.Cells.Item(88, 10).Value = 2300.8
$ws.Cells.Item(88, 11).Value = 2667.5
$ws.Cells.Item(88, 12).Value = 2300.8
$ws.Cells.Item(88, 13).Value = -2261.5
$ws.Cells.Item(88, 14).Value = -3112.8
$ws.Cells.Item(89, 8).Value = 8919.583000000001
$ws.Cells.Item(89, 9).Value = 3708
$ws.Cells.Item(89, 10).Value = 12642.143
$ws.Cells.Item(89, 11).Value = 18540
$ws.Cells.Item(89, 12).Value = 63210.715
$ws.Cells.Item(89, 13).Value = -12924
$ws.Cells.Item(89, 14).Value = -74442.715
$ws.Cells.Item(91, 8).Value = 2438.3125
$ws.Cells.Item(91, 9).Value = 2667.5
$ws.Cells.Item(91, 10).Value = 2300.8
$ws.Cells.Item(91, 11).Value = 2667.5
$ws.Cells.Item(91, 12).Value = 2300.8
$ws.Cells.Item(91, 13).Value = -1263.5
$ws.Cells.Item(91, 14).Value = -5108.8
$ws.Cells.Item(103, 8).Value = 23810424
$ws.Cells.Item(103, 9).Value = 761.6667
$ws.Cells.Item(103, 11).Value = 2285.0001
$ws.Cells.Item(103, 13).Value = -1699.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2951.818
$ws.Cells.Item(32, 9).Value = 2582.817
$ws.Cells.Item(32, 11).Value = 2582.817
$ws.Cells.Item(32, 13).Value = -2295.817
$ws.Cells.Item(41, 8).Value = 12606.2
$ws.Cells.Item(41, 9).Value = 13757.75
$ws.Cells.Item(41, 11).Value = 13757.75
$ws.Cells.Item(41, 13).Value = -13343.75
$ws.Cells.Item(74, 8).Value = 2436.0588
$ws.Cells.Item(74, 9).Value = 2496.9167
$ws.Cells.Item(74, 11).Value = 2496.9167
$ws.Cells.Item(74, 13).Value = -1622.9167
$ws.Cells.Item(77, 8).Value = 2436.0588
$ws.Cells.Item(77, 9).Value = 2496.9167
$ws.Cells.Item(77, 11).Value = 12484.5835
$ws.Cells.Item(77, 13).Value = -8116.583500000001
$ws.Cells.Item(80, 8).Value = 139983.33
$ws.Cells.Item(80, 10).Value = 159975
$ws.Cells.Item(80, 12).Value = 159975
$ws.Cells.Item(80, 14).Value = -161971
$ws.Cells.Item(83, 8).Value = 139983.33
$ws.Cells.Item(83, 10).Value = 159975
$ws.Cells.Item(83, 12).Value = 479925
$ws.Cells.Item(83, 14).Value = -489909
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 13).Value = ""
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 13).Value = ""
$ws.Cells.Item(97, 8).Value = 2029.2307
$ws.Cells.Item(97, 9).Value = 1943.6364
$ws.Cells.Item(97, 10).Value = 2500
$ws.Cells.Item(97, 11).Value = 1943.6364
$ws.Cells.Item(97, 12).Value = 2500
$ws.Cells.Item(97, 13).Value = -1447.6364
$ws.Cells.Item(97, 14).Value = -3492
$ws.Cells.Item(122, 8).Value = 2633.4634
$ws.Cells.Item(122, 9).Value = 2675.7
$ws.Cells.Item(122, 10).Value = 944
$ws.Cells.Item(122, 11).Value = 8027.099999999999
$ws.Cells.Item(122, 12).Value = 2832
$ws.Cells.Item(122, 13).Value = -5577.099999999999
$ws.Cells.Item(122, 14).Value = -7732
$ws.Cells.Item(132, 8).Value = 7815563.5
$ws.Cells.Item(132, 9).Value = 3118.7334
$ws.Cells.Item(132, 10).Value = 125002240
$ws.Cells.Item(132, 11).Value = 9356.200199999999
$ws.Cells.Item(132, 12).Value = 375006720
$ws.Cells.Item(132, 13).Value = -6826.200199999999
$ws.Cells.Item(132, 14).Value = -375011780

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 13).Value = ""
$ws.Cells.Item(105, 8).Value = 620996.3
$ws.Cells.Item(105, 9).Value = 1272610.1
$ws.Cells.Item(105, 10).Value = 3677.9473
$ws.Cells.Item(105, 11).Value = 1272610.1
$ws.Cells.Item(105, 12).Value = 3677.9473
$ws.Cells.Item(105, 13).Value = -1270863.1
$ws.Cells.Item(105, 14).Value = -7171.9473
$ws.Cells.Item(134, 8).Value = 2498.5217
$ws.Cells.Item(134, 9).Value = 2393
$ws.Cells.Item(134, 10).Value = 2999.75
$ws.Cells.Item(134, 11).Value = 7179
$ws.Cells.Item(134, 12).Value = 8999.25
$ws.Cells.Item(134, 13).Value = -4644
$ws.Cells.Item(134, 14).Value = -14069.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 62749.5
$ws.Cells.Item(28, 10).Value = 78666
$ws.Cells.Item(28, 12).Value = 78666
$ws.Cells.Item(28, 14).Value = -79156
$ws.Cells.Item(31, 8).Value = 1888.25
$ws.Cells.Item(31, 9).Value = 1680.4108
$ws.Cells.Item(31, 11).Value = 1680.4108
$ws.Cells.Item(31, 13).Value = -1385.4108
$ws.Cells.Item(34, 8).Value = 1888.25
$ws.Cells.Item(34, 9).Value = 1680.4108
$ws.Cells.Item(34, 11).Value = 1680.4108
$ws.Cells.Item(34, 13).Value = -1478.4108
$ws.Cells.Item(92, 8).Value = 35000
$ws.Cells.Item(92, 10).Value = 35000
$ws.Cells.Item(92, 12).Value = 35000
$ws.Cells.Item(92, 14).Value = -39992
$ws.Cells.Item(93, 8).Value = 59816.5
$ws.Cells.Item(93, 9).Value = 47225.25
$ws.Cells.Item(93, 11).Value = 47225.25
$ws.Cells.Item(93, 13).Value = -45353.25
$ws.Cells.Item(95, 8).Value = 19917.428
$ws.Cells.Item(95, 10).Value = 19917.428
$ws.Cells.Item(95, 12).Value = 19917.428
$ws.Cells.Item(95, 14).Value = -25409.428
$ws.Cells.Item(97, 8).Value = 40399
$ws.Cells.Item(97, 10).Value = 40399
$ws.Cells.Item(97, 12).Value = 40399
$ws.Cells.Item(97, 14).Value = -42381

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 1003.6667
$ws.Cells.Item(8, 9).Value = 1003.6667
$ws.Cells.Item(8, 11).Value = 3011.0001
$ws.Cells.Item(8, 13).Value = -2872.0001
$ws.Cells.Item(55, 8).Value = 5948.6665
$ws.Cells.Item(55, 9).Value = 2457.1428
$ws.Cells.Item(55, 10).Value = 18169
$ws.Cells.Item(55, 11).Value = 7371.428400000001
$ws.Cells.Item(55, 12).Value = 54507
$ws.Cells.Item(55, 13).Value = -7194.428400000001
$ws.Cells.Item(55, 14).Value = -54861
$ws.Cells.Item(57, 8).Value = 13828.25
$ws.Cells.Item(57, 9).Value = 5990
$ws.Cells.Item(57, 10).Value = 21666.5
$ws.Cells.Item(57, 11).Value = 17970
$ws.Cells.Item(57, 12).Value = 64999.5
$ws.Cells.Item(57, 13).Value = -17411
$ws.Cells.Item(57, 14).Value = -66117.5
$ws.Cells.Item(107, 8).Value = 4339901.5
$ws.Cells.Item(107, 9).Value = 5562.5
$ws.Cells.Item(107, 10).Value = 5359746
$ws.Cells.Item(107, 11).Value = 16687.5
$ws.Cells.Item(107, 12).Value = 16079238
$ws.Cells.Item(107, 13).Value = -14767.5
$ws.Cells.Item(107, 14).Value = -16083078
$ws.Cells.Item(138, 8).Value = 15544.357
$ws.Cells.Item(138, 9).Value = 19069.715
$ws.Cells.Item(138, 10).Value = 12019
$ws.Cells.Item(138, 11).Value = 57209.145
$ws.Cells.Item(138, 12).Value = 36057
$ws.Cells.Item(138, 13).Value = -52069.145
$ws.Cells.Item(138, 14).Value = -46337
$ws.Cells.Item(141, 8).Value = 7496.7144
$ws.Cells.Item(141, 9).Value = 3190.6667
$ws.Cells.Item(141, 11).Value = 9572.000100000001
$ws.Cells.Item(141, 13).Value = -4392.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4373.7407
$ws.Cells.Item(132, 9).Value = 4338.4
$ws.Cells.Item(132, 10).Value = 4417.9165
$ws.Cells.Item(132, 11).Value = 13015.2
$ws.Cells.Item(132, 12).Value = 13253.7495
$ws.Cells.Item(132, 13).Value = -10485.2
$ws.Cells.Item(132, 14).Value = -18313.7495

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 17861422
$ws.Cells.Item(100, 9).Value = 3491.5
$ws.Cells.Item(100, 10).Value = 62506250
$ws.Cells.Item(100, 11).Value = 3491.5
$ws.Cells.Item(100, 12).Value = 62506250
$ws.Cells.Item(100, 13).Value = -2950.5
$ws.Cells.Item(100, 14).Value = -62507332
$ws.Cells.Item(115, 8).Value = 0
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 14).Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 9224.583000000001
$ws.Cells.Item(62, 9).Value = 6066.6665
$ws.Cells.Item(62, 10).Value = 10277.223
$ws.Cells.Item(62, 11).Value = 6066.6665
$ws.Cells.Item(62, 12).Value = 10277.223
$ws.Cells.Item(62, 13).Value = -5442.6665
$ws.Cells.Item(62, 14).Value = -11525.223
$ws.Cells.Item(65, 8).Value = 9224.583000000001
$ws.Cells.Item(65, 9).Value = 6066.6665
$ws.Cells.Item(65, 10).Value = 10277.223
$ws.Cells.Item(65, 11).Value = 30333.3325
$ws.Cells.Item(65, 12).Value = 51386.115
$ws.Cells.Item(65, 13).Value = -27213.3325
$ws.Cells.Item(65, 14).Value = -57626.115
$ws.Cells.Item(81, 8).Value = 1999.8182
$ws.Cells.Item(81, 9).Value = 1444.2222
$ws.Cells.Item(81, 11).Value = 2888.4444
$ws.Cells.Item(81, 13).Value = -1827.4444
$ws.Cells.Item(84, 8).Value = 1999.8182
$ws.Cells.Item(84, 9).Value = 1444.2222
$ws.Cells.Item(84, 11).Value = 14442.222
$ws.Cells.Item(84, 13).Value = -9138.222
$ws.Cells.Item(122, 8).Value = 2059.1904
$ws.Cells.Item(122, 9).Value = 1663.6
$ws.Cells.Item(122, 10).Value = 3048.1667
$ws.Cells.Item(122, 11).Value = 4990.799999999999
$ws.Cells.Item(122, 12).Value = 9144.500100000001
$ws.Cells.Item(122, 13).Value = -2540.799999999999
$ws.Cells.Item(122, 14).Value = -14044.5001
$ws.Cells.Item(132, 8).Value = 3057.111
$ws.Cells.Item(132, 9).Value = 2566.375
$ws.Cells.Item(132, 11).Value = 7699.125
$ws.Cells.Item(132, 13).Value = -5169.125
$ws.Cells.Item(141, 8).Value = 225765
$ws.Cells.Item(141, 10).Value = 241111.11
$ws.Cells.Item(141, 12).Value = 241111.11
$ws.Cells.Item(141, 14).Value = -251471.11
